$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 762.8889
$ws.Range("I19").Value = 441.125
$ws.Range("J19").Value = 1020.3
$ws.Range("K19").Value = 441.125
$ws.Range("L19").Value = 1020.3
$ws.Range("M19").Value = -266.125
$ws.Range("N19").Value = -1370.3
$ws.Range("H107").Value = 452.43332
$ws.Range("I107").Value = 361.04166
$ws.Range("J107").Value = 818
$ws.Range("K107").Value = 361.04166
$ws.Range("L107").Value = 818
$ws.Range("M107").Value = 1558.95834
$ws.Range("N107").Value = -4658
$ws.Range("H125").Value = 14755
$ws.Range("I125").Value = 499.75
$ws.Range("J125").Value = 33762
$ws.Range("K125").Value = 4497.75
$ws.Range("L125").Value = 303858
$ws.Range("M125").Value = -2037.75
$ws.Range("N125").Value = -308778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4170170.2
$ws.Range("I2").Value = 4170170.2
$ws.Range("K2").Value = 4170170.2
$ws.Range("M2").Value = -4170057.2
$ws.Range("H97").Value = 685.4545000000001
$ws.Range("I97").Value = 646.6667
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 646.6667
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = -150.6667
$ws.Range("N97").Value = -2492
$ws.Range("H116").Value = 4170170.2
$ws.Range("I116").Value = 4170170.2
$ws.Range("K116").Value = 4170170.2
$ws.Range("M116").Value = -4167876.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4170170.2
$ws.Range("I3").Value = 4170170.2
$ws.Range("K3").Value = 4170170.2
$ws.Range("M3").Value = -4170056.2
$ws.Range("H94").Value = 1268.2354
$ws.Range("I94").Value = 1320
$ws.Range("J94").Value = 1100
$ws.Range("K94").Value = 1320
$ws.Range("L94").Value = 1100
$ws.Range("M94").Value = -869
$ws.Range("N94").Value = -2002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2166151
$ws.Range("I58").Value = 2756020.5
$ws.Range("K58").Value = 2756020.5
$ws.Range("M58").Value = -2755817.5
$ws.Range("H105").Value = 827.5
$ws.Range("I105").Value = 773.3333
$ws.Range("K105").Value = 773.3333
$ws.Range("M105").Value = 973.6667
$ws.Range("H136").Value = 2166151
$ws.Range("I136").Value = 2756020.5
$ws.Range("K136").Value = 8268061.5
$ws.Range("M136").Value = -8265511.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2393.6553
$ws.Range("I140").Value = 1782.5
$ws.Range("J140").Value = 2825.0588
$ws.Range("K140").Value = 5347.5
$ws.Range("L140").Value = 8475.1764
$ws.Range("M140").Value = -167.5
$ws.Range("N140").Value = -18835.1764

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3846
$ws.Range("I102").Value = 3147.0667
$ws.Range("K102").Value = 3147.0667
$ws.Range("M102").Value = -1525.0667
$ws.Range("H113").Value = 2156.0557
$ws.Range("I113").Value = 2134.0833
$ws.Range("J113").Value = 2200
$ws.Range("K113").Value = 2134.0833
$ws.Range("L113").Value = 2200
$ws.Range("M113").Value = 35.91670000000022
$ws.Range("N113").Value = -6540
$ws.Range("H122").Value = 6456.4546
$ws.Range("I122").Value = 8759
$ws.Range("J122").Value = 2427
$ws.Range("K122").Value = 26277
$ws.Range("L122").Value = 7281
$ws.Range("M122").Value = -23827
$ws.Range("N122").Value = -12181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3189.9092
$ws.Range("J7").Value = 2547.7
$ws.Range("L7").Value = 2547.7
$ws.Range("N7").Value = -2771.7
$ws.Range("H16").Value = 1240.6
$ws.Range("J16").Value = 1901
$ws.Range("L16").Value = 1901
$ws.Range("N16").Value = -2241
$ws.Range("H40").Value = 2861.2727
$ws.Range("I40").Value = 2492.4167
$ws.Range("J40").Value = 3303.9
$ws.Range("K40").Value = 2492.4167
$ws.Range("L40").Value = 3303.9
$ws.Range("M40").Value = -2356.4167
$ws.Range("N40").Value = -3575.9
$ws.Range("H68").Value = 2550
$ws.Range("I68").Value = 2100
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 2100
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -1351
$ws.Range("N68").Value = -4498
$ws.Range("H71").Value = 2550
$ws.Range("I71").Value = 2100
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 10500
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -6756
$ws.Range("N71").Value = -22488
$ws.Range("H100").Value = 2744.087
$ws.Range("I100").Value = 2190.2104
$ws.Range("J100").Value = 5375
$ws.Range("K100").Value = 2190.2104
$ws.Range("L100").Value = 5375
$ws.Range("M100").Value = -1649.2104
$ws.Range("N100").Value = -6457
$ws.Range("H122").Value = 6261.489
$ws.Range("I122").Value = 5853.6
$ws.Range("J122").Value = 7077.2666
$ws.Range("K122").Value = 17560.8
$ws.Range("L122").Value = 21231.7998
$ws.Range("M122").Value = -15110.8
$ws.Range("N122").Value = -26131.7998
$ws.Range("H126").Value = 3189.9092
$ws.Range("J126").Value = 2547.7
$ws.Range("L126").Value = 7643.099999999999
$ws.Range("N126").Value = -12583.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 11767324
$ws.Range("I81").Value = 1490
$ws.Range("J81").Value = 28575658
$ws.Range("K81").Value = 2980
$ws.Range("L81").Value = 57151316
$ws.Range("M81").Value = -1919
$ws.Range("N81").Value = -57153438
$ws.Range("H84").Value = 11767324
$ws.Range("I84").Value = 1490
$ws.Range("J84").Value = 28575658
$ws.Range("K84").Value = 14900
$ws.Range("L84").Value = 285756580
$ws.Range("M84").Value = -9596
$ws.Range("N84").Value = -285767188
$ws.Range("H113").Value = 591.5143
$ws.Range("I113").Value = 399.5
$ws.Range("J113").Value = 719.5238000000001
$ws.Range("K113").Value = 1198.5
$ws.Range("L113").Value = 2158.5714
$ws.Range("M113").Value = 971.5
$ws.Range("N113").Value = -6498.571400000001
$ws.Range("H126").Value = 1291
$ws.Range("I126").Value = 1253.6842
$ws.Range("K126").Value = 3761.0526
$ws.Range("M126").Value = -1291.0526
$ws.Range("H136").Value = 5470.5454
$ws.Range("I136").Value = 2387.7307
$ws.Range("J136").Value = 9923.5
$ws.Range("K136").Value = 7163.1921
$ws.Range("L136").Value = 29770.5
$ws.Range("M136").Value = -4613.1921
$ws.Range("N136").Value = -34870.5
